$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (rows 190 and 191)
$rows = @(
    @{ Row = 190; A = 45506.2916666667; B = 0;    C = 2.85999989509583; D = 2.85999989509583; E = 2.85999989509583; F = 2.85999989509583; G = "2.85999989509583"; H = "XHS.MI" },
    @{ Row = 191; A = 45509.6030671296; B = 2000; C = 2.79999995231628; D = 2.6800000667572;  E = 2.83999991416931; F = 2.79999995231628; G = "2.79999995231628"; H = "XHS.MI" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Column A (date) re-uses the same date/time style as every other row
    # in the column. Copy the formatting from the row above rather than
    # setting NumberFormat directly so we reuse the existing style index
    # instead of minting a brand-new (duplicate) one.
    $prevA = $ws.Cells.Item($rowNum - 1, 1)
    $cellA = $ws.Cells.Item($rowNum, 1)
    $prevA.Copy($cellA)
    $cellA.Value = $r.A

    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F

    # Column G (adj_close) is stored as TEXT in the source workbook, even
    # though it looks numeric. Force text type, assign, then drop the
    # number format back to the default "Normal" style so no style index
    # is left on the cell (matches the source formatting).
    $cellG = $ws.Cells.Item($rowNum, 7)
    $cellG.NumberFormat = "@"
    $cellG.Value = $r.G
    $cellG.Style = "Normal"

    $ws.Cells.Item($rowNum, 8).Value = $r.H
}
